# Add new power plant technologies to the Electricity Source subscript
# (issues #280 and #99) on the "ETS" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# New technology names to append, in order, starting at row 19.
$newTechs = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
$lastCol = 32   # column AF

for ($i = 0; $i -lt $newTechs.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTechs[$i]
    for ($col = 2; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# Put the selection on cell A25 for both sheets, matching the saved view state.
# Select the ETS sheet's cell first, then return to the About sheet so that it
# remains the active/selected tab, matching the original workbook state.
$ws.Range("A25").Select()
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Activate()
$aboutWs.Range("A25").Select()
